# FINFLUX-3612 Cartias specific scenarios
# Updates charge/penalty recalculation figures across the workbook
# (Summary, Repayment schedule, Transactions, ChargesTab) following the
# penalty amount change from $11.67 to $9.24, plus the related
# transaction-id renumbering on the Transactions sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 3466.18
$wsSummary.Range("E2").Value = 1533.82
$wsSummary.Range("A3").Value = 148.14
$wsSummary.Range("B3").Value = 124.58
$wsSummary.Range("E3").Value = 23.56
$wsSummary.Range("A5").Value = 9.24
$wsSummary.Range("B5").Value = 9.24
$wsSummary.Activate()
$wsSummary.Range("C9").Select()

# ---------------------------------------------------------------
# Repayment schedule
# ---------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("H5").Value = 27.08
$wsRepay.Range("J5").Value = 9.24
$wsRepay.Range("K5").Value = 896.96
$wsRepay.Range("N5").Value = 963.92
$wsRepay.Range("Q5").Value = 36.32
$wsRepay.Range("F6").Value = 742.89
$wsRepay.Range("G6").Value = 754.61
$wsRepay.Range("H6").Value = 15.12
$wsRepay.Range("F7").Value = 749.96
$wsRepay.Range("G7").Value = 4.65
$wsRepay.Range("H7").Value = 8.05
$wsRepay.Range("F8").Value = 4.65
$wsRepay.Range("H8").Value = 0.39
$wsRepay.Range("K8").Value = 5.04
$wsRepay.Range("Q8").Value = 5.04
$wsRepay.Activate()
$wsRepay.Range("J9").Select()

# ---------------------------------------------------------------
# Transactions
# ---------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")
$wsTxn.Range("A2").Value = 382
$wsTxn.Range("E2").Value = 41
$wsTxn.Range("J2").Value = 1531.81
$wsTxn.Range("A3").Value = 381
$wsTxn.Range("E3").Value = 41
$wsTxn.Range("A4").Value = 380
$wsTxn.Range("A5").Value = 377
$wsTxn.Range("A6").Value = 376
$wsTxn.Range("A7").Value = 375
$wsTxn.Range("A8").Value = 374
$wsTxn.Range("A9").Value = 373
$wsTxn.Activate()
$wsTxn.Range("G6").Select()

# ---------------------------------------------------------------
# ChargesTab - the $11.67 penalty text becomes $9.24
# ---------------------------------------------------------------
$wsCharges = $wb.Worksheets.Item("ChargesTab")
$wsCharges.Range("G3").Value = "$9.24"
$wsCharges.Range("H3").Value = "$9.24"
$wsCharges.Rows.Item(3).RowHeight = 45
$wsCharges.Columns.Item(6).ColumnWidth = 9.85546875
$wsCharges.Activate()
$wsCharges.Range("F9").Select()

# Leave Transactions as the active sheet, matching the workbook state.
$wsTxn.Activate()
